$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Toplam" (Total) header cell (F1) is being dropped from the export.
# Deleting it with a shift-left moves every following header one column to
# the left (G->F, H->G, I->H, J->I, K->J) and drops the now-unused "Toplam"
# shared string, shrinking the sheet's used range by one column.
$ws.Range("F1").Delete(-4159)

# Correct/relabel the currency columns (now I1 and J1 after the shift).
$ws.Range("I1").Value = "Para birimi"
$ws.Range("J1").Value = "Para Birimi Tutarı"

# Restore the saved selection/active-cell position.
$ws.Range("J2").Select() | Out-Null
